$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "29.151.69"
Set-TextValue "D3" "1.864.27"
$ws.Range("E3").Value = "  -1.02%  "
Set-TextValue "D4" "0.9998"
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "0.7101"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  -0.46%  "
Set-TextValue "D7" "1.000"
$ws.Range("E7").Value = "  +0.09%  "
Set-TextValue "D8" "0.3095"
$ws.Range("E8").Value = "  -0.89%  "
Set-TextValue "D9" "0.07645"
$ws.Range("E9").Value = "  -3.71%  "
Set-TextValue "D10" "24.63"
$ws.Range("E10").Value = "  -2.98%  "
Set-TextValue "D11" "0.08357"
$ws.Range("E11").Value = "  +0.80%  "
Set-TextValue "D12" "1.886.53"
$ws.Range("E12").Value = "  +0.14%  "
Set-TextValue "D13" "5.227"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("E14").Value = "  -3.30%  "
Set-TextValue "D15" "91.23"
$ws.Range("E15").Value = "  -0.41%  "
Set-TextValue "D16" "29.196.06"
$ws.Range("E16").Value = "  -1.00%  "
Set-TextValue "D17" "5.922"
$ws.Range("E17").Value = "  -0.65%  "
Set-TextValue "D18" "243.11"
$ws.Range("E18").Value = "  -2.06%  "
Set-TextValue "D19" "0.000007811"
$ws.Range("E19").Value = "  -0.96%  "
Set-TextValue "D20" "2.113.06"
$ws.Range("E20").Value = "  -0.78%  "
Set-TextValue "D21" "13.07"
$ws.Range("E21").Value = "  -2.35%  "
Set-TextValue "D22" "0.9997"
$ws.Range("E22").Value = "  +0.12%  "
Set-TextValue "D23" "7.859"
$ws.Range("E23").Value = "  -1.67%  "
Set-TextValue "D24" "1.000"
$ws.Range("E24").Value = "  +0.10%  "
Set-TextValue "D25" "0.1581"
$ws.Range("E25").Value = "  -2.32%  "
Set-TextValue "D26" "163.13"
$ws.Range("E26").Value = "  -0.17%  "
Set-TextValue "D27" "8.936"
$ws.Range("E27").Value = "  -1.57%  "
Set-TextValue "D28" "18.44"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  -0.19%  "
Set-TextValue "D30" "1.324"
$ws.Range("E30").Value = "  -2.50%  "
Set-TextValue "D31" "4.390"
$ws.Range("E31").Value = "  -0.39%  "
Set-TextValue "D32" "4.248"
Set-TextValue "D33" "0.05154"
$ws.Range("E33").Value = "  -2.72%  "
Set-TextValue "D34" "0.7946"
$ws.Range("E34").Value = "  +8.97%  "
Set-TextValue "D35" "1.910"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("E36").Value = "  -3.38%  "
Set-TextValue "D37" "2.687"
$ws.Range("E37").Value = "  +0.43%  "
Set-TextValue "D38" "0.01846"
Set-TextValue "D39" "2.707"
$ws.Range("E39").Value = "  -1.02%  "
Set-TextValue "D40" "1.163.92"
$ws.Range("E40").Value = "  -5.19%  "
Set-TextValue "D41" "6.249"
$ws.Range("E41").Value = "  +0.63%  "
Set-TextValue "D42" "0.8915"
$ws.Range("E42").Value = "  -2.38%  "
Set-TextValue "D43" "72.91"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  -0.16%  "
Set-TextValue "D46" "2.010.32"
$ws.Range("E46").Value = "  -1.27%  "
Set-TextValue "D47" "0.5180"
$ws.Range("E47").Value = "  -1.83%  "
Set-TextValue "D48" "1.772"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  +0.34%  "
